$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width change: column L (12) from 27 to 28 characters
$ws.Columns.Item(12).ColumnWidth = 27.14

# Cell value updates
$ws.Range("E2").Value = "2026-02-11 21:18:44"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "82%"
$ws.Range("I2").Value = "9.0 mm"
$ws.Range("E3").Value = "2026-02-11 21:18:46"
$ws.Range("L3").Value = "79.9 km/h - 214º 20:59 TU"
$ws.Range("M3").Value = "5.5 °C 20:48 TU"
$ws.Range("O3").Value = "0.4 °C"
$ws.Range("E4").Value = "2026-02-11 21:18:49"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "57%"
$ws.Range("J4").Value = "1002.0 hPa"
$ws.Range("E5").Value = "2026-02-11 21:18:52"
$ws.Range("M5").Value = "4.6 °C 20:34 TU"
$ws.Range("O5").Value = "0.8 °C"
$ws.Range("E6").Value = "2026-02-11 21:18:54"
$ws.Range("J6").Value = "1002.5 hPa"
$ws.Range("O6").Value = "13.3 °C"
$ws.Range("E7").Value = "2026-02-11 21:18:57"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "46%"
$ws.Range("J7").Value = "1003.3 hPa"
$ws.Range("O7").Value = "18.9 °C"
$ws.Range("E8").Value = "2026-02-11 21:19:00"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "58%"
$ws.Range("J8").Value = "1003.0 hPa"
$ws.Range("N8").Value = "10.9 °C 20:42 TU"
$ws.Range("O8").Value = "14.9 °C"
$ws.Range("E9").Value = "2026-02-11 21:19:02"
$ws.Range("I9").Value = "2.8 mm"
$ws.Range("E10").Value = "2026-02-11 21:19:05"
$ws.Range("O10").Value = "13.7 °C"
$ws.Range("E11").Value = "2026-02-11 21:19:08"
$ws.Range("I11").Value = "1.5 mm"
$ws.Range("O11").Value = "7.9 °C"
$ws.Range("E12").Value = "2026-02-11 21:19:10"
$ws.Range("E13").Value = "2026-02-11 21:19:13"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "80%"
$ws.Range("I13").Value = "2.2 mm"
$ws.Range("O13").Value = "7.3 °C"
$ws.Range("E14").Value = "2026-02-11 21:19:16"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "51%"
$ws.Range("N14").Value = "13.6 °C 20:58 TU"
$ws.Range("O14").Value = "18.6 °C"
$ws.Range("E15").Value = "2026-02-11 21:19:18"
$ws.Range("I15").Value = "3.1 mm"
$ws.Range("E16").Value = "2026-02-11 21:19:20"
$ws.Range("I16").Value = "8.8 mm"
$ws.Range("L16").Value = "107.3 km/h - 224º 20:50 TU"
$ws.Range("E17").Value = "2026-02-11 21:19:23"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "79%"
$ws.Range("E18").Value = "2026-02-11 21:19:26"
$ws.Range("J18").Value = "1002.6 hPa"
$ws.Range("E19").Value = "2026-02-11 21:19:28"
$ws.Range("E20").Value = "2026-02-11 21:19:31"
$ws.Range("I20").Value = "3.3 mm"
$ws.Range("O20").Value = "-0.9 °C"
$ws.Range("E21").Value = "2026-02-11 21:19:34"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "84%"
$ws.Range("I21").Value = "4.5 mm"
$ws.Range("J21").Value = "1005.4 hPa"
$ws.Range("E22").Value = "2026-02-11 21:19:36"
$ws.Range("L22").Value = "109.8 km/h - 297º 20:58 TU"
$ws.Range("E23").Value = "2026-02-11 21:19:39"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "75%"
$ws.Range("E24").Value = "2026-02-11 21:19:42"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "77%"
$ws.Range("J24").Value = "1006.5 hPa"
$ws.Range("N24").Value = "10.2 °C 20:47 TU"
$ws.Range("E25").Value = "2026-02-11 21:19:44"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "66%"
$ws.Range("I25").Value = "4.2 mm"
$ws.Range("E26").Value = "2026-02-11 21:19:47"
$ws.Range("I26").Value = "2.6 mm"
$ws.Range("J26").Value = "1002.4 hPa"
$ws.Range("E27").Value = "2026-02-11 21:19:50"
$ws.Range("I27").Value = "2.9 mm"
$ws.Range("E28").Value = "2026-02-11 21:19:52"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "78%"
$ws.Range("J28").Value = "1002.7 hPa"
$ws.Range("E29").Value = "2026-02-11 21:19:55"
$ws.Range("L29").Value = "59.0 km/h - 242º 20:31 TU"
$ws.Range("O29").Value = "13.2 °C"
$ws.Range("E30").Value = "2026-02-11 21:19:57"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "88%"
$ws.Range("I30").Value = "6.0 mm"
$ws.Range("J30").Value = "1002.6 hPa"
$ws.Range("K30").Value = "10.5 MJ/m2"
$ws.Range("E31").Value = "2026-02-11 21:20:00"
$ws.Range("I31").Value = "3.2 mm"
$ws.Range("J31").Value = "1001.8 hPa"
$ws.Range("E32").Value = "2026-02-11 21:20:03"
$ws.Range("E33").Value = "2026-02-11 21:20:05"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "82%"
$ws.Range("I33").Value = "3.3 mm"
$ws.Range("J33").Value = "1004.6 hPa"
$ws.Range("E34").Value = "2026-02-11 21:20:08"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "63%"
$ws.Range("I34").Value = "4.9 mm"
$ws.Range("E35").Value = "2026-02-11 21:20:11"
$ws.Range("E36").Value = "2026-02-11 21:20:14"
$ws.Range("I36").Value = "8.4 mm"
$ws.Range("J36").Value = "1002.8 hPa"
$ws.Range("O36").Value = "13.1 °C"
$ws.Range("E37").Value = "2026-02-11 21:20:16"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "81%"
$ws.Range("I37").Value = "1.4 mm"
$ws.Range("J37").Value = "1004.0 hPa"
$ws.Range("E38").Value = "2026-02-11 21:20:19"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "60%"
$ws.Range("E39").Value = "2026-02-11 21:20:22"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "58%"
$ws.Range("I39").Value = "3.6 mm"
$ws.Range("E40").Value = "2026-02-11 21:20:24"
$ws.Range("I40").Value = "6.9 mm"
$ws.Range("J40").Value = "1006.6 hPa"
$ws.Range("E41").Value = "2026-02-11 21:20:27"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "51%"
$ws.Range("J41").Value = "1004.4 hPa"
$ws.Range("O41").Value = "18.5 °C"
$ws.Range("E42").Value = "2026-02-11 21:20:30"
$ws.Range("O42").Value = "12.9 °C"
$ws.Range("E43").Value = "2026-02-11 21:20:32"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "66%"
$ws.Range("N43").Value = "9.7 °C 20:59 TU"
$ws.Range("O43").Value = "12.9 °C"
$ws.Range("E44").Value = "2026-02-11 21:20:35"
$ws.Range("I44").Value = "7.6 mm"
$ws.Range("L44").Value = "91.1 km/h - 204º 20:38 TU"
$ws.Range("M44").Value = "2.8 °C 20:41 TU"
$ws.Range("O44").Value = "-0.4 °C"
$ws.Range("E45").Value = "2026-02-11 21:20:37"
$ws.Range("J45").Value = "1005.2 hPa"
$ws.Range("O45").Value = "6.9 °C"
$ws.Range("E46").Value = "2026-02-11 21:20:40"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "63%"
$ws.Range("I46").Value = "2.0 mm"
$ws.Range("J46").Value = "1006.8 hPa"
$ws.Range("N46").Value = "12.2 °C 20:58 TU"
$ws.Range("O46").Value = "16.7 °C"
